# Monitored endpoint already integrated with message
# -----------------------------------------------------------------
# This script updates the "last_values" style report workbook:
#  - renames the sheet
#  - tweaks the indicator/label text
#  - refreshes the timestamp value + format
#  - removes now-unused helper cells (B3 / F5)
#  - formats the header row (bold, bordered, centered)
# -----------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (Planilha1 -> Sheet1)
$ws.Name = "Sheet1"

# --- Text / label updates -------------------------------------------------
# A1: "indicator" -> "info"
$ws.Range("A1").Value = "info"

# A2: "monitored_total" -> "monitored_goods"
$ws.Range("A2").Value = "monitored_goods"

# --- Timestamp update -------------------------------------------------
# B2 moves from 2023-06-01 (45078) to 2023-07-18 (45125) and gets a
# full date-time display format instead of the short date format.
# Reset to the plain "Normal" style first so no leftover alignment /
# border / bold formatting carries over from the old style.
$ws.Range("B2").Style = "Normal"
$ws.Range("B2").Value = 45125
$ws.Range("B2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- Remove now-unused cells -------------------------------------------------
# B3 was an empty, number-formatted placeholder cell - drop it entirely.
$ws.Range("B3").Clear()
# F5 was an empty, styled placeholder cell - drop it entirely.
$ws.Range("F5").Clear()

# --- Header row formatting -------------------------------------------------
# A1:B1 become a bold, centered, thin-bordered header. Reset to "Normal"
# first so the old column-level number format (#,##0.00) doesn't linger.
$header = $ws.Range("A1:B1")
$header.Style = "Normal"
$header.Font.Bold = $true
$header.Borders.LineStyle = 1
$header.HorizontalAlignment = -4108   # xlCenter
$header.VerticalAlignment = -4160     # xlTop

# Tidy up row heights back to the default (no more custom 19.5pt rows)
$ws.Rows.Item(1).RowHeight = 15
$ws.Rows.Item(2).RowHeight = 15
$ws.Rows.Item(3).RowHeight = 15
$ws.Rows.Item(4).RowHeight = 15
$ws.Rows.Item(5).RowHeight = 15

Write-Output "Workbook updated: header/values refreshed, unused cells removed."
